# Update Odds.xlsx ("updated odds with poe.how")
# Applies the Dropchance/Talented/85-Area value updates from the commit,
# plus the resulting view-state (scroll position / selection) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dropchance (C), Talented (D), 85 Area (E) updates -------------------

# shaper
$ws.Range("C4").Value = 51.5
$ws.Range("C5").Value = 24
$ws.Range("C6").Value = 16
$ws.Range("C7").Value = 3.5
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 10.47
$ws.Range("C12").Value = 3
$ws.Range("E12").Value = 3

# eater of worlds
$ws.Range("C23").Value = 45
$ws.Range("C24").Value = 33
$ws.Range("C25").Value = 5
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 2
$ws.Range("E29").Value = 35

# searing exarch
$ws.Range("C30").Value = 33
$ws.Range("C31").Value = 33
$ws.Range("C32").Value = 10
$ws.Range("C34").Value = 8
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 2
$ws.Range("E36").Value = 35

# atziri
$ws.Range("C37").Value = 67.97
$ws.Range("C38").Value = 12.2
$ws.Range("C39").Value = 4.58
$ws.Range("C42").Value = 15.25
$ws.Range("C43").Value = 15.47
$ws.Range("C44").Value = 10.68
$ws.Range("C45").Value = 5.23
$ws.Range("C46").Value = 10.02
$ws.Range("C47").Value = 0

# uber elder
$ws.Range("C50").Value = 16.88
$ws.Range("C51").Value = 0
$ws.Range("E51").Value = 14
$ws.Range("C52").Value = 14
$ws.Range("C53").Value = 7.5
$ws.Range("C54").Value = 3
$ws.Range("E54").Value = 2

# uber atziri
$ws.Range("C62").Value = 3.48
$ws.Range("C65").Value = 46.52
$ws.Range("C67").Value = 37.81

# maven
$ws.Range("C71").Value = 26.15
$ws.Range("C72").Value = 16.92
$ws.Range("C73").Value = 24.62
$ws.Range("C74").Value = 16.92
$ws.Range("C75").Value = 12.31
$ws.Range("C77").Value = 2.44
$ws.Range("C79").Value = 2.36
$ws.Range("C80").Value = 49.23
$ws.Range("E80").Formula = "=100-C80"

# sirus
$ws.Range("C84").Value = 34.72

# --- View state: scrolled down a few rows, selection moved to C92 --------

$excel.Goto($ws.Range("A52"), $true)
$ws.Range("C92").Select()
